$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 0.1794871794871795
$ws.Cells.Item(2, 3).Value = 0.5811965811965812
$ws.Cells.Item(2, 10).Value = 0.0170940170940171
$ws.Cells.Item(2, 16).Value = 0.1367521367521368
$ws.Cells.Item(2, 19).Value = 0.08547008547008547
$ws.Cells.Item(3, 2).Value = 0.01449275362318841
$ws.Cells.Item(3, 3).Value = 0.01449275362318841
$ws.Cells.Item(3, 10).Value = 0.007246376811594203
$ws.Cells.Item(3, 16).Value = 0.7246376811594203
$ws.Cells.Item(3, 19).Value = 0.2391304347826087
$ws.Cells.Item(4, 10).Value = 0.02325581395348837
$ws.Cells.Item(4, 16).Value = 0.6976744186046512
$ws.Cells.Item(4, 19).Value = 0.2790697674418605
$ws.Cells.Item(6, 2).Value = 0.06363636363636363
$ws.Cells.Item(6, 4).Value = 0.01818181818181818
$ws.Cells.Item(6, 6).Value = 0.06818181818181818
$ws.Cells.Item(6, 10).Value = 0.25
$ws.Cells.Item(6, 15).Value = 0.01818181818181818
$ws.Cells.Item(6, 17).Value = 0.1545454545454545
$ws.Cells.Item(6, 18).Value = 0.1045454545454545
$ws.Cells.Item(6, 19).Value = 0.3227272727272728
$ws.Cells.Item(7, 2).Value = 0.1047120418848168
$ws.Cells.Item(7, 4).Value = 0.01570680628272251
$ws.Cells.Item(7, 5).Value = 0.005235602094240838
$ws.Cells.Item(7, 6).Value = 0.04712041884816754
$ws.Cells.Item(7, 10).Value = 0.1518324607329843
$ws.Cells.Item(7, 15).Value = 0.005235602094240838
$ws.Cells.Item(7, 17).Value = 0.1884816753926702
$ws.Cells.Item(7, 18).Value = 0.0418848167539267
$ws.Cells.Item(7, 19).Value = 0.4397905759162304
$ws.Cells.Item(8, 2).Value = 0.07093821510297482
$ws.Cells.Item(8, 4).Value = 0.02517162471395881
$ws.Cells.Item(8, 6).Value = 0.06864988558352403
$ws.Cells.Item(8, 10).Value = 0.1052631578947368
$ws.Cells.Item(8, 15).Value = 0.02059496567505721
$ws.Cells.Item(8, 17).Value = 0.1739130434782609
$ws.Cells.Item(8, 18).Value = 0.09610983981693363
$ws.Cells.Item(8, 19).Value = 0.4393592677345537
$ws.Cells.Item(9, 2).Value = 0.09580838323353294
$ws.Cells.Item(9, 4).Value = 0.01796407185628742
$ws.Cells.Item(9, 6).Value = 0.04191616766467066
$ws.Cells.Item(9, 10).Value = 0.08982035928143713
$ws.Cells.Item(9, 15).Value = 0.02994011976047904
$ws.Cells.Item(9, 17).Value = 0.1616766467065868
$ws.Cells.Item(9, 18).Value = 0.08982035928143713
$ws.Cells.Item(9, 19).Value = 0.4730538922155689
$ws.Cells.Item(10, 2).Value = 0.08713692946058091
$ws.Cells.Item(10, 4).Value = 0.02157676348547718
$ws.Cells.Item(10, 5).Value = 0.0008298755186721991
$ws.Cells.Item(10, 6).Value = 0.08132780082987552
$ws.Cells.Item(10, 10).Value = 0.1037344398340249
$ws.Cells.Item(10, 15).Value = 0.01493775933609959
$ws.Cells.Item(10, 17).Value = 0.2232365145228216
$ws.Cells.Item(10, 18).Value = 0.08298755186721991
$ws.Cells.Item(10, 19).Value = 0.3842323651452282
$ws.Cells.Item(11, 7).Value = 0.1242424242424242
$ws.Cells.Item(11, 10).Value = 0.1060606060606061
$ws.Cells.Item(11, 11).Value = 0.203030303030303
$ws.Cells.Item(11, 12).Value = 0.5606060606060606
$ws.Cells.Item(11, 19).Value = 0.006060606060606061
$ws.Cells.Item(12, 7).Value = 0.7010309278350515
$ws.Cells.Item(12, 10).Value = 0.1855670103092784
$ws.Cells.Item(12, 11).Value = 0.005154639175257732
$ws.Cells.Item(12, 12).Value = 0.03608247422680412
$ws.Cells.Item(12, 19).Value = 0.07216494845360824
$ws.Cells.Item(15, 6).Value = 0.01619433198380567
$ws.Cells.Item(15, 8).Value = 0.1781376518218623
$ws.Cells.Item(15, 9).Value = 0.09716599190283401
$ws.Cells.Item(15, 10).Value = 0.3319838056680162
$ws.Cells.Item(15, 11).Value = 0.08502024291497975
$ws.Cells.Item(15, 13).Value = 0.008097165991902834
$ws.Cells.Item(15, 15).Value = 0.07692307692307693
$ws.Cells.Item(15, 19).Value = 0.2064777327935223
$ws.Cells.Item(16, 6).Value = 0.01923076923076923
$ws.Cells.Item(16, 8).Value = 0.1666666666666667
$ws.Cells.Item(16, 9).Value = 0.0641025641025641
$ws.Cells.Item(16, 10).Value = 0.3782051282051282
$ws.Cells.Item(16, 11).Value = 0.1666666666666667
$ws.Cells.Item(16, 13).Value = 0.00641025641025641
$ws.Cells.Item(16, 14).Value = 0.00641025641025641
$ws.Cells.Item(16, 15).Value = 0.08333333333333333
$ws.Cells.Item(16, 19).Value = 0.108974358974359
$ws.Cells.Item(17, 6).Value = 0.0272108843537415
$ws.Cells.Item(17, 8).Value = 0.2063492063492063
$ws.Cells.Item(17, 9).Value = 0.08616780045351474
$ws.Cells.Item(17, 10).Value = 0.3786848072562358
$ws.Cells.Item(17, 11).Value = 0.08616780045351474
$ws.Cells.Item(17, 13).Value = 0.01360544217687075
$ws.Cells.Item(17, 15).Value = 0.07256235827664399
$ws.Cells.Item(17, 19).Value = 0.1292517006802721
$ws.Cells.Item(18, 6).Value = 0.01058201058201058
$ws.Cells.Item(18, 8).Value = 0.1851851851851852
$ws.Cells.Item(18, 9).Value = 0.07407407407407407
$ws.Cells.Item(18, 10).Value = 0.4444444444444444
$ws.Cells.Item(18, 11).Value = 0.09523809523809523
$ws.Cells.Item(18, 13).Value = 0.01058201058201058
$ws.Cells.Item(18, 15).Value = 0.05291005291005291
$ws.Cells.Item(18, 19).Value = 0.126984126984127
$ws.Cells.Item(19, 6).Value = 0.01826846703733121
$ws.Cells.Item(19, 8).Value = 0.193010325655282
$ws.Cells.Item(19, 9).Value = 0.06433677521842732
$ws.Cells.Item(19, 10).Value = 0.3717235901509134
$ws.Cells.Item(19, 11).Value = 0.1239078633836378
$ws.Cells.Item(19, 13).Value = 0.01906274821286736
$ws.Cells.Item(19, 14).Value = 0.00238284352660842
$ws.Cells.Item(19, 15).Value = 0.08101667990468626
$ws.Cells.Item(19, 19).Value = 0.1262907069102462

Write-Output "Updated team specific time data values for Illinois St._B sheet."
